$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (A2m | Lrp1 | FAPs -> ECs) - new TPM values
$ws.Range("M2").Value = 3.456265333333333
$ws.Range("N2").Value = 10.368796
$ws.Range("O2").Value = 0.009841535807677501
$ws.Range("P2").Value = 0.0098415358076775
$ws.Range("Q2").Value = 0.6195470818844444
$ws.Range("R2").Value = 5.57592373696
$ws.Range("S2").Value = 0.009841535807677501
$ws.Range("T2").Value = 0.0098415358076775

# Row 3 (A2m | Lrp1 | FAPs -> FAPs) - derived specificities recomputed
$ws.Range("O3").Value = 0.8587907398420774
$ws.Range("P3").Value = 0.8587907398420773
$ws.Range("S3").Value = 0.8587907398420774
$ws.Range("T3").Value = 0.8587907398420773

# Row 4 (A2m | Lrp1 | FAPs -> MuSCs) - derived specificities recomputed
$ws.Range("O4").Value = 0.1313677243502452
$ws.Range("P4").Value = 0.1313677243502452
$ws.Range("S4").Value = 0.1313677243502452
$ws.Range("T4").Value = 0.1313677243502452
